$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.504.77"
$ws.Cells.Item(2, 5).Value = "  +2.31%  "
$ws.Cells.Item(3, 4).Value = "1.873.82"
$ws.Cells.Item(3, 5).Value = "  +1.71%  "
$ws.Cells.Item(4, 4).Value = "1.015"
$ws.Cells.Item(4, 5).Value = "  +1.01%  "
$ws.Cells.Item(5, 4).Value = "313.42"
$ws.Cells.Item(5, 5).Value = "  +1.39%  "
$ws.Cells.Item(6, 5).Value = "  +0.92%  "
$ws.Cells.Item(7, 4).Value = "0.4784"
$ws.Cells.Item(7, 5).Value = "  +0.84%  "
$ws.Cells.Item(8, 4).Value = "0.3772"
$ws.Cells.Item(8, 5).Value = "  +2.62%  "
$ws.Cells.Item(9, 4).Value = "0.07380"
$ws.Cells.Item(9, 5).Value = "  +2.48%  "
$ws.Cells.Item(10, 4).Value = "0.9379"
$ws.Cells.Item(10, 5).Value = "  +1.30%  "
$ws.Cells.Item(11, 4).Value = "20.73"
$ws.Cells.Item(11, 5).Value = "  +5.64%  "
$ws.Cells.Item(12, 4).Value = "0.07852"
$ws.Cells.Item(12, 5).Value = "  +2.78%  "
$ws.Cells.Item(13, 4).Value = "1.905.10"
$ws.Cells.Item(13, 5).Value = "  +2.92%  "
$ws.Cells.Item(14, 4).Value = "5.449"
$ws.Cells.Item(14, 5).Value = "  +2.62%  "
$ws.Cells.Item(15, 4).Value = "6.590"
$ws.Cells.Item(15, 5).Value = "  +2.86%  "
$ws.Cells.Item(16, 4).Value = "90.93"
$ws.Cells.Item(16, 5).Value = "  +2.48%  "
$ws.Cells.Item(17, 4).Value = "1.015"
$ws.Cells.Item(17, 5).Value = "  +0.72%  "
$ws.Cells.Item(18, 4).Value = "0.000008939"
$ws.Cells.Item(18, 5).Value = "  +3.44%  "
$ws.Cells.Item(19, 4).Value = "1.013"
$ws.Cells.Item(19, 5).Value = "  +0.88%  "
$ws.Cells.Item(20, 5).Value = "  +2.62%  "
$ws.Cells.Item(21, 4).Value = "27.556.44"
$ws.Cells.Item(21, 5).Value = "  +2.41%  "
$ws.Cells.Item(22, 4).Value = "5.138"
$ws.Cells.Item(22, 5).Value = "  +1.85%  "
$ws.Cells.Item(23, 5).Value = "  +0.88%  "
$ws.Cells.Item(24, 5).Value = "  +2.15%  "
$ws.Cells.Item(25, 4).Value = "153.99"
$ws.Cells.Item(25, 5).Value = "  +1.22%  "
$ws.Cells.Item(26, 4).Value = "18.55"
$ws.Cells.Item(26, 5).Value = "  +2.32%  "
$ws.Cells.Item(27, 4).Value = "2.021"
$ws.Cells.Item(27, 5).Value = "  +1.02%  "
$ws.Cells.Item(28, 4).Value = "116.03"
$ws.Cells.Item(28, 5).Value = "  +1.56%  "
$ws.Cells.Item(29, 4).Value = "5.002"
$ws.Cells.Item(29, 5).Value = "  +1.11%  "
$ws.Cells.Item(30, 4).Value = "0.08938"
$ws.Cells.Item(30, 5).Value = "  +0.99%  "
$ws.Cells.Item(31, 4).Value = "3.342"
$ws.Cells.Item(31, 5).Value = "  +1.33%  "
$ws.Cells.Item(32, 4).Value = "1.218"
$ws.Cells.Item(32, 5).Value = "  +3.48%  "
$ws.Cells.Item(33, 4).Value = "4.623"
$ws.Cells.Item(33, 5).Value = "  +3.04%  "
$ws.Cells.Item(34, 4).Value = "0.7524"
$ws.Cells.Item(34, 5).Value = "  +0.29%  "
$ws.Cells.Item(35, 4).Value = "2.691"
$ws.Cells.Item(35, 5).Value = "  -2.14%  "
$ws.Cells.Item(36, 4).Value = "0.02069"
$ws.Cells.Item(36, 5).Value = "  +6.25%  "
$ws.Cells.Item(37, 4).Value = "1.117"
$ws.Cells.Item(37, 5).Value = "  +2.59%  "
$ws.Cells.Item(38, 4).Value = "0.05306"
$ws.Cells.Item(38, 5).Value = "  +0.80%  "
$ws.Cells.Item(39, 4).Value = "3.005"
$ws.Cells.Item(39, 5).Value = "  +1.50%  "
$ws.Cells.Item(40, 4).Value = "0.5349"
$ws.Cells.Item(40, 5).Value = "  +2.37%  "
$ws.Cells.Item(41, 4).Value = "7.095"
$ws.Cells.Item(41, 5).Value = "  +2.34%  "
$ws.Cells.Item(42, 5).Value = "  +0.89%  "
$ws.Cells.Item(43, 4).Value = "8.428"
$ws.Cells.Item(43, 5).Value = "  +2.50%  "
$ws.Cells.Item(44, 4).Value = "10.60"
$ws.Cells.Item(44, 5).Value = "  +0.36%  "
$ws.Cells.Item(45, 4).Value = "0.4837"
$ws.Cells.Item(45, 5).Value = "  +2.27%  "
$ws.Cells.Item(46, 5).Value = "  +0.91%  "
$ws.Cells.Item(47, 4).Value = "1.663"
$ws.Cells.Item(47, 5).Value = "  +3.57%  "
$ws.Cells.Item(48, 4).Value = "103.18"
$ws.Cells.Item(48, 5).Value = "  +1.47%  "
$ws.Cells.Item(49, 5).Value = "  +2.83%  "
$ws.Cells.Item(50, 4).Value = "0.06095"
$ws.Cells.Item(50, 5).Value = "  +1.18%  "
$ws.Cells.Item(51, 4).Value = "0.8977"
$ws.Cells.Item(51, 5).Value = "  +1.44%  "
